# Apply the "add 2022-Q4 data" edit:
#  1. Insert a new worksheet named "2022-Q4" right before the "2022-Q3"
#     sheet and populate it with the quarterly fund-holdings table (same
#     header/column layout and styling as the other quarter sheets). The
#     sheet is produced by duplicating "2022-Q3" (so header row, the 0-based
#     index column and the boxed header/index style all come along for
#     free) and then overwriting the data columns.
#  2. Update the "总计" (summary) sheet: add a 2022-Q4 row at the top of the
#     data and push the remaining quarters down by one row (column A keeps
#     its running 0-based index; a new row 7 is appended for 2021-Q3).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. New "2022-Q4" sheet
# ---------------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3)
$q4 = $wb.Worksheets.Item("2022-Q3 (2)")
$q4.Name = "2022-Q4"

# "2022-Q3" has 12 data rows (rows 2-13); "2022-Q4" only needs 9 (rows
# 2-10), so drop the extra rows copied along with the sheet.
$q4.Rows("11:13").Delete()

# Fund code / size / position columns look numeric but are stored as text
# (leading zeros, fixed 2/4-decimal strings) - force text format before
# writing so the leading zeros / trailing zeros survive.
$q4.Range("B2:B10").NumberFormat = "@"
$q4.Range("D2:G10").NumberFormat = "@"

$rows = @(
    @("012526", "广发盛锦混合A",              "24.00", "93.13", "4.02", "0.9648", 8),
    @("002851", "南方品质优选灵活配置混合A",    "11.52", "85.40", "6.49", "0.7476", 5),
    @("005123", "南方优享分红灵活配置混合A",    "6.79",  "92.15", "6.91", "0.4692", 7),
    @("006587", "南方优享分红灵活配置混合C",    "1.84",  "92.15", "6.91", "0.1271", 7),
    @("001692", "南方国策动力股票",            "3.58",  "87.60", "1.88", "0.0673", 9),
    @("012527", "广发盛锦混合C",              "1.14",  "93.13", "4.02", "0.0458", 8),
    @("012426", "南方价值臻选混合A",            "3.85",  "75.34", "0.59", "0.0227", 10),
    @("012427", "南方价值臻选混合C",            "0.19",  "75.34", "0.59", "0.0011", 10),
    @("013501", "南方品质优选灵活配置混合C",    "0.01",  "85.40", "6.49", "0.0006", 5)
)

$r = 2
foreach ($row in $rows) {
    $q4.Cells.Item($r, 2).Value = $row[0]
    $q4.Cells.Item($r, 3).Value = $row[1]
    $q4.Cells.Item($r, 4).Value = $row[2]
    $q4.Cells.Item($r, 5).Value = $row[3]
    $q4.Cells.Item($r, 6).Value = $row[4]
    $q4.Cells.Item($r, 7).Value = $row[5]
    $q4.Cells.Item($r, 8).Value = $row[6]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# 2. "总计" summary sheet
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Append row 7 (2021-Q3, previously the last row) - clone the index cell's
# style from A6, then fill in the values for the appended quarter.
$total.Cells.Item(6, 1).Copy($total.Cells.Item(7, 1))
$total.Cells.Item(7, 1).Value = 5
$total.Cells.Item(7, 2).Value = "2021-Q3"
$total.Cells.Item(7, 3).Value = 6
$total.Cells.Item(7, 4).Value = 1.49

# Shift the quarter figures down one slot so the newest quarter lands on
# top (column A's running index is untouched - it already reads 0..4).
$total.Cells.Item(6, 2).Value = "2021-Q4"
$total.Cells.Item(6, 3).Value = 3
$total.Cells.Item(6, 4).Value = 0.54

$total.Cells.Item(5, 2).Value = "2022-Q1"
$total.Cells.Item(5, 3).Value = 10
$total.Cells.Item(5, 4).Value = 1.07

$total.Cells.Item(4, 2).Value = "2022-Q2"
$total.Cells.Item(4, 3).Value = 9
$total.Cells.Item(4, 4).Value = 0.66

$total.Cells.Item(3, 2).Value = "2022-Q3"
$total.Cells.Item(3, 3).Value = 12
$total.Cells.Item(3, 4).Value = 1.41

$total.Cells.Item(2, 2).Value = "2022-Q4"
$total.Cells.Item(2, 3).Value = 9
$total.Cells.Item(2, 4).Value = 2.45
